$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows 66 and 67 ---
$ws.Range("J66").Value = 12289217
$ws.Range("M66").Value = 915803
$ws.Range("Q66").Value = -221646
$ws.Range("W66").Value = 13577169
$ws.Range("X66").Value = -1507792

$ws.Range("J67").Value = 15919356
$ws.Range("M67").Value = 93763
$ws.Range("Q67").Value = -2379102
$ws.Range("W67").Value = 17735819
$ws.Range("X67").Value = -4194398

# --- New row 68 (01-07-2021 quarter) ---
# Entering "01-07-2021" directly gets auto-detected as a date serial, which
# is not what the source data wants (column A stores plain text labels like
# the other quarters, e.g. "01-04-2021"). Building the text via a formula
# sidesteps the date auto-detection, then Copy/PasteSpecial(values) bakes it
# back down to a plain literal (shared string) without leaving any leftover
# cell formatting/style behind.
$ws.Range("A68").Formula = '="01-07-2021"'
$ws.Range("A68").Copy()
$ws.Range("A68").PasteSpecial(-4163)

$ws.Range("B68").Value = 14115745
$ws.Range("C68").Value = 11243045
$ws.Range("D68").Value = 855625
$ws.Range("E68").Value = 699242
$ws.Range("F68").Value = 18295
$ws.Range("G68").Value = 166543
$ws.Range("H68").Value = 255784
$ws.Range("I68").Value = 877211
$ws.Range("J68").Value = 20209553
$ws.Range("K68").Value = 2858534
$ws.Range("L68").Value = 1113924
$ws.Range("M68").Value = 881840
$ws.Range("N68").Value = 13292928
$ws.Range("O68").Value = 2021446
$ws.Range("P68").Value = 40881
$ws.Range("Q68").Value = -6093807
$ws.Range("R68").Value = 1822745
$ws.Range("S68").Value = 1326
$ws.Range("T68").Value = 976024
$ws.Range("U68").Value = 848047
$ws.Range("V68").Value = 14117071
$ws.Range("W68").Value = 22033623
$ws.Range("X68").Value = -7916552
$ws.Range("Y68").Value = -3
